$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "StartGameCommandMessage"
$ws.Range("B12").Value = 1030
$ws.Range("C12").Value = "Notifica i client che è stata avviata una partita"

$ws.Range("A12").Select()
